$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The label "glycemie moyenne estimee" used to live in column A (row 8) as
# part of the list of excluded-variable names. It is being moved into
# column B instead (row 10, as an additionally-included variable), which
# removes it from column A and shifts the remaining column-A entries
# (rows 9:60) up by one row, leaving column A one row shorter overall.
# Column B is left untouched aside from the new value at B10.

# Capture the column-A values that need to move up one row (old rows 9:60
# become new rows 8:59), then write them back and clear the now-empty
# trailing row.
$lastRow = 60
$colAValues = @()
for ($r = 9; $r -le $lastRow; $r++) {
    $colAValues += , ($ws.Cells.Item($r, 1).Value2)
}
for ($i = 0; $i -lt $colAValues.Length; $i++) {
    $ws.Cells.Item(8 + $i, 1).Value2 = $colAValues[$i]
}
$ws.Cells.Item($lastRow, 1).ClearContents()

# Place the relocated label into column B.
$ws.Range("B10").Value2 = "glycemie moyenne estimee"

# Update the visible selection to match the new view of the sheet.
$ws.Range("C68").Select()
